# Adds the latest weekly price records for "Cebollín" (Vega Monumental
# Concepción) by inserting two new rows right above the existing data
# block (at row 72), pushing the previous rows down by two, and filling
# the two freshly inserted rows with the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 72-73; everything currently at 72:89 shifts
# down to 74:91 (dates/prices for those rows stay exactly as they were).
$ws.Range("A72:R73").EntireRow.Insert()

# New row 72 - "Primera" quality, Región de Ñuble, $/paquete 6 unidades
$ws.Range("A72").Value = 11
$ws.Range("B72").Value = "Vega Monumental Concepción"
$ws.Range("C72").Value = "Bíobío"
$ws.Range("D72").Value2 = 44951
$ws.Range("E72").Value = 8
$ws.Range("F72").Value = 100112037
$ws.Range("G72").Value = "Cebollín"
$ws.Range("H72").Value = "Sin especificar"
$ws.Range("I72").Value = "Primera"
$ws.Range("J72").Value = 200
$ws.Range("K72").Value = 700
$ws.Range("L72").Value = 800
$ws.Range("M72").Value = 750
$ws.Range("N72").Value = "`$/paquete 6 unidades"
$ws.Range("O72").Value = "Región de Ñuble"
$ws.Range("P72").Value = 125
$ws.Range("Q72").Value = 6
$ws.Range("R72").Value = "Hortaliza"

# New row 73 - "Segunda" quality, Región de Ñuble, $/paquete 6 unidades
$ws.Range("A73").Value = 11
$ws.Range("B73").Value = "Vega Monumental Concepción"
$ws.Range("C73").Value = "Bíobío"
$ws.Range("D73").Value2 = 44951
$ws.Range("E73").Value = 8
$ws.Range("F73").Value = 100112037
$ws.Range("G73").Value = "Cebollín"
$ws.Range("H73").Value = "Sin especificar"
$ws.Range("I73").Value = "Segunda"
$ws.Range("J73").Value = 100
$ws.Range("K73").Value = 600
$ws.Range("L73").Value = 600
$ws.Range("M73").Value = 600
$ws.Range("N73").Value = "`$/paquete 6 unidades"
$ws.Range("O73").Value = "Región de Ñuble"
$ws.Range("P73").Value = 100
$ws.Range("Q73").Value = 6
$ws.Range("R73").Value = "Hortaliza"
